$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 193 ("Vega Modelo de Temuco" /
# Granada price-series). Excel shifts rows 193:264 down to 194:265 and grows
# the sheet's used range to A1:T265.
$ws.Rows.Item(193).Insert()

# Populate the newly inserted row 193 with the new observation.
$ws.Cells.Item(193, 1).Value = 10
$ws.Cells.Item(193, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(193, 3).Value = "La Araucanía"
$ws.Cells.Item(193, 4).Value = 45120
$ws.Cells.Item(193, 5).Value = 9
$ws.Cells.Item(193, 6).Value = "Fruta"
$ws.Cells.Item(193, 7).Value = 100104
$ws.Cells.Item(193, 8).Value = "Frutos de pepita"
$ws.Cells.Item(193, 9).Value = 100104001
$ws.Cells.Item(193, 10).Value = "Granada"
$ws.Cells.Item(193, 11).Value = "Wonderfull"
$ws.Cells.Item(193, 12).Value = "Primera"
$ws.Cells.Item(193, 13).Value = 300
$ws.Cells.Item(193, 14).Value = 12000
$ws.Cells.Item(193, 15).Value = 13000
$ws.Cells.Item(193, 16).Value = 12333
$ws.Cells.Item(193, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(193, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(193, 19).Value = 1233
$ws.Cells.Item(193, 20).Value = 10
